# OPR344_ExportManifest_TestData.xlsx edit
# Commit message: "added  trait field and lte testcases"
#
# 1. Add three new worksheets (OPR344_EXP_00003/4/5) with LTE (lying-list /
#    priority) export-manifest test rows, the last of which has a new
#    "NewFlightNumber" trait column (sheet4 gets "SplitPieces" trait column).
# 2. Sheet2 loses its tabSelected flag (new sheet5 becomes the selected tab).
# 3. A couple of cells hold a quote-prefixed numeric-looking flight number
#    ("0316") that must stay text, entered with a leading apostrophe.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 3: OPR344_EXP_00003
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$ws3.Name = "OPR344_EXP_00003"

$ws3.Range("A1").Value = "AgentCode"
$ws3.Range("B1").Value = "ShipperCode "
$ws3.Range("C1").Value = "ConsigneeCode"
$ws3.Range("D1").Value = "Origin"
$ws3.Range("E1").Value = "Destination"
$ws3.Range("F1").Value = "ProductCode"
$ws3.Range("G1").Value = "SCC"
$ws3.Range("H1").Value = "Commodity"
$ws3.Range("I1").Value = "ShipmentDescription"
$ws3.Range("J1").Value = "ServiceCargoClass"
$ws3.Range("K1").Value = "Piece"
$ws3.Range("L1").Value = "Weight"
$ws3.Range("M1").Value = "ChargeType"
$ws3.Range("N1").Value = "ModeOfPayment"
$ws3.Range("O1").Value = "cartType"
$ws3.Range("P1").Value = "AWBSectionName"

$ws3.Range("A2").Value = 11377
$ws3.Range("B2").Value = 11377
$ws3.Range("C2").Value = 11377
$ws3.Range("D2").Value = "SEA"
$ws3.Range("E2").Value = "LAX"
$ws3.Range("F2").Value = "GENERAL "
$ws3.Range("G2").Value = "None"
$ws3.Range("H2").Value = "NONSCR"
$ws3.Range("I2").Value = "None"
$ws3.Range("J2").Value = "None"
$ws3.Range("K2").Value = 13
$ws3.Range("L2").Value = 775
$ws3.Range("M2").Value = "CC"
$ws3.Range("N2").Value = "None"
$ws3.Range("O2").Value = "CART"
$ws3.Range("P2").Value = "PlannedShipment"

$ws3.Range("A3").Value = 11377
$ws3.Range("B3").Value = 11377
$ws3.Range("C3").Value = 11377
$ws3.Range("D3").Value = "ANC"
$ws3.Range("E3").Value = "HNL"
$ws3.Range("F3").Value = "PRIORITY "
$ws3.Range("G3").Value = "None"
$ws3.Range("H3").Value = 2199
$ws3.Range("I3").Value = "None"
$ws3.Range("J3").Value = "None"
$ws3.Range("K3").Value = 8
$ws3.Range("L3").Value = 360
$ws3.Range("M3").Value = "CC"
$ws3.Range("N3").Value = "None"
$ws3.Range("O3").Value = "CART"
$ws3.Range("P3").Value = "PlannedShipment"

$ws3.Range("A4").Value = 11377
$ws3.Range("B4").Value = 11377
$ws3.Range("C4").Value = 11377
$ws3.Range("D4").Value = "SAN"
$ws3.Range("E4").Value = "JFK"
$ws3.Range("F4").Value = "GOLDSTREAK"
$ws3.Range("G4").Value = "None"
$ws3.Range("H4").Value = "NONSCR"
$ws3.Range("I4").Value = "None"
$ws3.Range("J4").Value = "None"
$ws3.Range("K4").Value = 2
$ws3.Range("L4").Value = 59
$ws3.Range("M4").Value = "CC"
$ws3.Range("N4").Value = "None"
$ws3.Range("O4").Value = "CART"
$ws3.Range("P4").Value = "PlannedShipment"

$ws3.Range("Q2").Select() | Out-Null

# ---------------------------------------------------------------------
# Sheet 4: OPR344_EXP_00004
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$ws4.Name = "OPR344_EXP_00004"

$ws4.Range("A1").Value = "AgentCode"
$ws4.Range("B1").Value = "ShipperCode "
$ws4.Range("C1").Value = "ConsigneeCode"
$ws4.Range("D1").Value = "Origin"
$ws4.Range("E1").Value = "Destination"
$ws4.Range("F1").Value = "ProductCode"
$ws4.Range("G1").Value = "SCC"
$ws4.Range("H1").Value = "Commodity"
$ws4.Range("I1").Value = "ShipmentDescription"
$ws4.Range("J1").Value = "ServiceCargoClass"
$ws4.Range("K1").Value = "Piece"
$ws4.Range("L1").Value = "Weight"
$ws4.Range("M1").Value = "ChargeType"
$ws4.Range("N1").Value = "ModeOfPayment"
$ws4.Range("O1").Value = "AWBSectionName"
$ws4.Range("P1").Value = "SplitPieces"
$ws4.Range("Q1").Value = "cartType"

$ws4.Range("A2").Value = 11377
$ws4.Range("B2").Value = 11377
$ws4.Range("C2").Value = 11377
$ws4.Range("D2").Value = "SAN"
$ws4.Range("E2").Value = "JFK"
$ws4.Range("F2").Value = "GENERAL"
$ws4.Range("G2").Value = "None"
$ws4.Range("H2").Value = "'0316"
$ws4.Range("I2").Value = "None"
$ws4.Range("J2").Value = "None"
$ws4.Range("K2").Value = 2
$ws4.Range("L2").Value = 234
$ws4.Range("M2").Value = "CC"
$ws4.Range("N2").Value = "None"
$ws4.Range("O2").Value = "PlannedShipment"
$ws4.Range("P2").Value = 1
$ws4.Range("Q2").Value = "CART"

$ws4.Range("A3").Value = 11377
$ws4.Range("B3").Value = 11377
$ws4.Range("C3").Value = 11377
$ws4.Range("D3").Value = "ANC"
$ws4.Range("E3").Value = "HNL"
$ws4.Range("F3").Value = "PRIORITY"
$ws4.Range("G3").Value = "None"
$ws4.Range("H3").Value = 2199
$ws4.Range("I3").Value = "None"
$ws4.Range("J3").Value = "None"
$ws4.Range("K3").Value = 10
$ws4.Range("L3").Value = 189
$ws4.Range("M3").Value = "CC"
$ws4.Range("N3").Value = "None"
$ws4.Range("O3").Value = "PlannedShipment"
$ws4.Range("P3").Value = 1
$ws4.Range("Q3").Value = "CART"

$ws4.Range("M7").Select() | Out-Null

# ---------------------------------------------------------------------
# Sheet 5: OPR344_EXP_00005
# ---------------------------------------------------------------------
$ws5 = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$ws5.Name = "OPR344_EXP_00005"

$ws5.Range("A1").Value = "AgentCode"
$ws5.Range("B1").Value = "ShipperCode "
$ws5.Range("C1").Value = "ConsigneeCode"
$ws5.Range("D1").Value = "Origin"
$ws5.Range("E1").Value = "Destination"
$ws5.Range("F1").Value = "ProductCode"
$ws5.Range("G1").Value = "SCC"
$ws5.Range("H1").Value = "Commodity"
$ws5.Range("I1").Value = "ShipmentDescription"
$ws5.Range("J1").Value = "ServiceCargoClass"
$ws5.Range("K1").Value = "Piece"
$ws5.Range("L1").Value = "Weight"
$ws5.Range("M1").Value = "ChargeType"
$ws5.Range("N1").Value = "ModeOfPayment"
$ws5.Range("O1").Value = "AWBSectionName"
$ws5.Range("P1").Value = "NewFlightNumber"
$ws5.Range("Q1").Value = "cartType"

$ws5.Range("A2").Value = 11377
$ws5.Range("B2").Value = 11377
$ws5.Range("C2").Value = 11377
$ws5.Range("D2").Value = "SEA"
$ws5.Range("E2").Value = "JFK"
$ws5.Range("F2").Value = "GENERAL"
$ws5.Range("G2").Value = "None"
$ws5.Range("H2").Value = "'0316"
$ws5.Range("I2").Value = "None"
$ws5.Range("J2").Value = "None"
$ws5.Range("K2").Value = 2
$ws5.Range("L2").Value = 59
$ws5.Range("M2").Value = "CC"
$ws5.Range("N2").Value = "None"
$ws5.Range("O2").Value = "PlannedShipment"
$ws5.Range("P2").Value = 26
$ws5.Range("Q2").Value = "CART"

$ws5.Range("C4").Select() | Out-Null

# Sheet5 becomes the active / selected tab (mirrors activeTab=4 in the
# saved workbook, and clears sheet2's tabSelected flag automatically).
$ws5.Activate()
